$wb = $excel.ActiveWorkbook

# --- time_variants sheet: remove the IPT-era "1950" and "1990" year columns ---
$tv = $wb.Worksheets.Item("time_variants")

# Delete the higher-lettered column first so the lower one's letter doesn't shift
# out from under us before we get to it.
$tv.Columns("G:G").Delete()
$tv.Columns("E:E").Delete()

# Scroll/selection state: make the sheet reselect near the top-left and put the
# cursor on A7 (was B7) in the bottom-right frozen pane.
$tv.Activate()
[void]$tv.Range("A7").Select()

# --- Active tab moves from time_variants back to constants ---
$constants = $wb.Worksheets.Item("constants")
$constants.Activate()
